$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell value updates: [cellRef, newValue, forceText]
$updates = @(
    @('D2', '29.420.47', $False),
    @('E2', '  -0.03%  ', $False),
    @('D3', '1.853.46', $False),
    @('E3', '  +0.24%  ', $False),
    @('D4', '0.9992', $True),
    @('E4', '  -0.09%  ', $False),
    @('D5', '240.91', $True),
    @('E5', '  +0.03%  ', $False),
    @('E6', '  +0.68%  ', $False),
    @('D7', '1.001', $True),
    @('E7', '  +0.05%  ', $False),
    @('D8', '0.07590', $True),
    @('E8', '  -1.44%  ', $False),
    @('D9', '0.2932', $True),
    @('E9', '  +0.35%  ', $False),
    @('D10', '24.54', $True),
    @('E10', '  -0.91%  ', $False),
    @('D11', '0.07758', $True),
    @('E11', '  +0.19%  ', $False),
    @('D12', '1.854.19', $False),
    @('E12', '  +0.37%  ', $False),
    @('D13', '5.032', $True),
    @('E13', '  +0.03%  ', $False),
    @('D14', '0.6818', $True),
    @('E14', '  +0.24%  ', $False),
    @('D15', '0.00001049', $True),
    @('E15', '  -1.95%  ', $False),
    @('D16', '83.39', $True),
    @('E16', '  -0.24%  ', $False),
    @('D17', '6.143', $True),
    @('E17', '  -0.58%  ', $False),
    @('D18', '29.414.89', $False),
    @('E18', '  -0.13%  ', $False),
    @('D19', '230.43', $True),
    @('E19', '  +1.01%  ', $False),
    @('D20', '12.38', $True),
    @('E20', '  -0.40%  ', $False),
    @('D21', '1.000', $True),
    @('E21', '  -0.02%  ', $False),
    @('D22', '7.491', $True),
    @('E22', '  +0.92%  ', $False),
    @('D23', '1.000', $True),
    @('E23', '  -0.08%  ', $False),
    @('D24', '159.23', $True),
    @('E24', '  +0.95%  ', $False),
    @('E25', '  +1.09%  ', $False),
    @('D26', '8.468', $True),
    @('E26', '  +0.77%  ', $False),
    @('D27', '17.69', $True),
    @('E27', '  -0.01%  ', $False),
    @('E28', '  +4.00%  ', $False),
    @('D29', '1.478', $True),
    @('E29', '  +0.72%  ', $False),
    @('D30', '0.05689', $True),
    @('E30', '  +0.18%  ', $False),
    @('D31', '4.128', $True),
    @('E31', '  +0.17%  ', $False),
    @('D32', '4.058', $True),
    @('E32', '  +0.72%  ', $False),
    @('D33', '1.833', $True),
    @('E33', '  -0.53%  ', $False),
    @('D34', '1.158', $True),
    @('E34', '  -0.34%  ', $False),
    @('D35', '0.7093', $True),
    @('E35', '  +0.19%  ', $False),
    @('D36', '2.583', $True),
    @('E36', '  -0.05%  ', $False),
    @('D37', '0.01830', $True),
    @('E37', '  +2.35%  ', $False),
    @('D38', '1.246.18', $False),
    @('E38', '  +2.07%  ', $False),
    @('E39', '  -1.80%  ', $False),
    @('D40', '6.440', $True),
    @('E40', '  -1.61%  ', $False),
    @('E41', '  -0.17%  ', $False),
    @('D42', '1.000', $True),
    @('E42', '  -0.07%  ', $False),
    @('D43', '2.010.33', $False),
    @('E43', '  -0.14%  ', $False),
    @('D44', '102.42', $True),
    @('E44', '  +0.66%  ', $False),
    @('D45', '65.96', $True),
    @('E45', '  -0.36%  ', $False),
    @('D46', '7.149', $True),
    @('E46', '  +0.08%  ', $False),
    @('D47', '0.1168', $True),
    @('E47', '  +1.81%  ', $False),
    @('B48', 'BabyDogeCoin', $False),
    @('C48', 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge', $False),
    @('D48', '0.00000000116', $True),
    @('E48', '  -4.92%  ', $False),
    @('B49', 'EnergySwap', $False),
    @('C49', 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens', $False),
    @('D49', '9.012', $True),
    @('E49', '  -0.23%  ', $False),
    @('D50', '0.3965', $True),
    @('E50', '  -1.41%  ', $False),
    @('D51', '1.676', $True),
    @('E51', '  +0.14%  ', $False)
)

foreach ($u in $updates) {
    $cellRef = $u[0]
    $val = $u[1]
    $forceText = $u[2]
    $rng = $ws.Range($cellRef)
    if ($forceText) {
        $rng.NumberFormat = "@"
        $rng.Value = $val
        $rng.ClearFormats()
    } else {
        $rng.Value = $val
    }
}